$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1366.6666
$ws.Range("I12").Value = 1366.6666
$ws.Range("K12").Value = 1366.6666
$ws.Range("M12").Value = -1196.6666
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("H40").Value = 1369.66
$ws.Range("I40").Value = 1173
$ws.Range("J40").Value = 1490.1936
$ws.Range("K40").Value = 1173
$ws.Range("L40").Value = 1490.1936
$ws.Range("M40").Value = -998
$ws.Range("N40").Value = -1840.1936
$ws.Range("H62").Value = 9615.280000000001
$ws.Range("I62").Value = 9053.799999999999
$ws.Range("K62").Value = 9053.799999999999
$ws.Range("M62").Value = -8429.799999999999
$ws.Range("H65").Value = 9615.280000000001
$ws.Range("I65").Value = 9053.799999999999
$ws.Range("K65").Value = 45269
$ws.Range("M65").Value = -42149
$ws.Range("H69").Value = 11916.167
$ws.Range("J69").Value = 12539.4
$ws.Range("L69").Value = 37618.2
$ws.Range("N69").Value = -39366.2
$ws.Range("H72").Value = 11916.167
$ws.Range("J72").Value = 12539.4
$ws.Range("L72").Value = 112854.6
$ws.Range("N72").Value = -121590.6
$ws.Range("H132").Value = 6079729.5
$ws.Range("I132").Value = 7430525.5
$ws.Range("J132").Value = 1149.8
$ws.Range("K132").Value = 22291576.5
$ws.Range("L132").Value = 3449.4
$ws.Range("M132").Value = -22289046.5
$ws.Range("N132").Value = -8509.4
$ws.Range("H137").Value = 1992.5333
$ws.Range("J137").Value = 1366.3334
$ws.Range("L137").Value = 4099.0002
$ws.Range("N137").Value = -9199.0002
$ws.Range("H141").Value = 995
$ws.Range("I141").Value = 995
$ws.Range("K141").Value = 2985
$ws.Range("M141").Value = 2195

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3873.9333
$ws.Range("I32").Value = 3132.3157
$ws.Range("J32").Value = 7899.857
$ws.Range("K32").Value = 3132.3157
$ws.Range("L32").Value = 7899.857
$ws.Range("M32").Value = -2845.3157
$ws.Range("N32").Value = -8473.857
$ws.Range("H61").Value = 4194.4
$ws.Range("I61").Value = 3496.7273
$ws.Range("K61").Value = 3496.7273
$ws.Range("M61").Value = -3284.7273
$ws.Range("H63").Value = 3112.5
$ws.Range("H66").Value = 3112.5
$ws.Range("H97").Value = 3318.2144
$ws.Range("I97").Value = 2573.4614
$ws.Range("K97").Value = 2573.4614
$ws.Range("M97").Value = -2077.4614
$ws.Range("H122").Value = 73530.07000000001
$ws.Range("I122").Value = 201319.8
$ws.Range("J122").Value = 2535.7778
$ws.Range("K122").Value = 603959.3999999999
$ws.Range("L122").Value = 7607.3334
$ws.Range("M122").Value = -601509.3999999999
$ws.Range("N122").Value = -12507.3334
$ws.Range("H132").Value = 2327.5925
$ws.Range("I132").Value = 1923.762
$ws.Range("K132").Value = 5771.286
$ws.Range("M132").Value = -3241.286
$ws.Range("H136").Value = 4194.4
$ws.Range("I136").Value = 3496.7273
$ws.Range("K136").Value = 10490.1819
$ws.Range("M136").Value = -7940.1819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7054.769
$ws.Range("I86").Value = 5390.3687
$ws.Range("K86").Value = 5390.3687
$ws.Range("M86").Value = -4267.3687
$ws.Range("H89").Value = 7054.769
$ws.Range("I89").Value = 5390.3687
$ws.Range("K89").Value = 26951.8435
$ws.Range("M89").Value = -21335.8435
$ws.Range("H94").Value = 1385.4445
$ws.Range("I94").Value = 922.6
$ws.Range("K94").Value = 922.6
$ws.Range("M94").Value = -471.6
$ws.Range("H99").Value = 5755.8237
$ws.Range("I99").Value = 2112.375
$ws.Range("J99").Value = 8994.444
$ws.Range("K99").Value = 2112.375
$ws.Range("L99").Value = 8994.444
$ws.Range("M99").Value = -614.375
$ws.Range("N99").Value = -11990.444
$ws.Range("H105").Value = 5011.875
$ws.Range("I105").Value = 5950
$ws.Range("J105").Value = 2197.5
$ws.Range("K105").Value = 5950
$ws.Range("L105").Value = 2197.5
$ws.Range("M105").Value = -4203
$ws.Range("N105").Value = -5691.5
$ws.Range("H134").Value = 2746.24
$ws.Range("I134").Value = 2176.6316
$ws.Range("K134").Value = 6529.8948
$ws.Range("M134").Value = -3994.8948

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 641.5
$ws.Range("J22").Value = 812.25
$ws.Range("L22").Value = 812.25
$ws.Range("N22").Value = -1512.25
$ws.Range("H58").Value = 3299.8333
$ws.Range("I58").Value = 3299.8333
$ws.Range("K58").Value = 3299.8333
$ws.Range("M58").Value = -3096.8333
$ws.Range("H80").Value = 40062
$ws.Range("J80").Value = 40062
$ws.Range("L80").Value = 40062
$ws.Range("N80").Value = -42308
$ws.Range("H83").Value = 40062
$ws.Range("J83").Value = 40062
$ws.Range("L83").Value = 120186
$ws.Range("N83").Value = -131418
$ws.Range("H132").Value = 3352.157
$ws.Range("I132").Value = 3286.9167
$ws.Range("K132").Value = 9860.750100000001
$ws.Range("M132").Value = -7330.750100000001
$ws.Range("H134").Value = 15687.8
$ws.Range("I134").Value = 9166.929
$ws.Range("K134").Value = 27500.787
$ws.Range("M134").Value = -24965.787
$ws.Range("H136").Value = 3299.8333
$ws.Range("I136").Value = 3299.8333
$ws.Range("K136").Value = 9899.499899999999
$ws.Range("M136").Value = -7349.499899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 604.4286
$ws.Range("I86").Value = 477.66666
$ws.Range("K86").Value = 1432.99998
$ws.Range("M86").Value = -246.9999800000001
$ws.Range("H89").Value = 604.4286
$ws.Range("I89").Value = 477.66666
$ws.Range("K89").Value = 4298.99994
$ws.Range("M89").Value = 1629.00006
$ws.Range("H140").Value = 3807.4546
$ws.Range("I140").Value = 3807.4546
$ws.Range("K140").Value = 11422.3638
$ws.Range("M140").Value = -6242.363799999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3686.8386
$ws.Range("I80").Value = 2812.1333
$ws.Range("J80").Value = 4506.875
$ws.Range("K80").Value = 2812.1333
$ws.Range("L80").Value = 4506.875
$ws.Range("M80").Value = -1814.1333
$ws.Range("N80").Value = -6502.875
$ws.Range("H83").Value = 3686.8386
$ws.Range("I83").Value = 2812.1333
$ws.Range("J83").Value = 4506.875
$ws.Range("K83").Value = 14060.6665
$ws.Range("L83").Value = 22534.375
$ws.Range("M83").Value = -9068.666499999999
$ws.Range("N83").Value = -32518.375
$ws.Range("H102").Value = 2940.7878
$ws.Range("I102").Value = 2881.3572
$ws.Range("K102").Value = 2881.3572
$ws.Range("M102").Value = -1259.3572
$ws.Range("H132").Value = 479644.44
$ws.Range("I132").Value = 628470.9
$ws.Range("K132").Value = 1885412.7
$ws.Range("M132").Value = -1882882.7

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 20839282
$ws.Range("I16").Value = 83335704
$ws.Range("J16").Value = 7141.9443
$ws.Range("K16").Value = 83335704
$ws.Range("L16").Value = 7141.9443
$ws.Range("M16").Value = -83335534
$ws.Range("N16").Value = -7481.9443
$ws.Range("H22").Value = 2575
$ws.Range("I22").Value = 2800
$ws.Range("J22").Value = 2350
$ws.Range("K22").Value = 2800
$ws.Range("L22").Value = 2350
$ws.Range("M22").Value = -2505
$ws.Range("N22").Value = -2940
$ws.Range("H27").Value = 2575
$ws.Range("I27").Value = 2800
$ws.Range("J27").Value = 2350
$ws.Range("K27").Value = 2800
$ws.Range("L27").Value = 2350
$ws.Range("M27").Value = -2693
$ws.Range("N27").Value = -2564
$ws.Range("H82").Value = 1764
$ws.Range("I82").Value = 1866.5
$ws.Range("J82").Value = 1600
$ws.Range("K82").Value = 1866.5
$ws.Range("L82").Value = 1600
$ws.Range("M82").Value = -1505.5
$ws.Range("N82").Value = -2322
$ws.Range("H85").Value = 1764
$ws.Range("I85").Value = 1866.5
$ws.Range("J85").Value = 1600
$ws.Range("K85").Value = 1866.5
$ws.Range("L85").Value = 1600
$ws.Range("M85").Value = -618.5
$ws.Range("N85").Value = -4096
$ws.Range("H120").Value = 30634
$ws.Range("I120").Value = 30634
$ws.Range("K120").Value = 30634
$ws.Range("M120").Value = -25796
$ws.Range("H122").Value = 4126.478
$ws.Range("I122").Value = 3612.5293
$ws.Range("K122").Value = 10837.5879
$ws.Range("M122").Value = -8387.5879

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 35035.145
$ws.Range("J63").Value = 35035.145
$ws.Range("L63").Value = 35035.145
$ws.Range("N63").Value = -36283.145
$ws.Range("H66").Value = 35035.145
$ws.Range("J66").Value = 35035.145
$ws.Range("L66").Value = 105105.435
$ws.Range("N66").Value = -111345.435
$ws.Range("H92").Value = 36500
$ws.Range("J92").Value = 36500
$ws.Range("L92").Value = 36500
$ws.Range("N92").Value = -41492
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H122").Value = 2296.3333
$ws.Range("I122").Value = 2086.389
$ws.Range("K122").Value = 6259.167
$ws.Range("M122").Value = -3809.167
